$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, thin border, centered) from H1
# onto the two new header cells so I1/J1 match the look of the other
# header cells (B1..H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

# New data values for row 3
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5
